# Auto-generated script to update cryptos price table
# Applies 117 cell text updates (columns B, C, D, E) across rows 2-51
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell value as literal text, preserving its original style
# (needed because plain numeric-looking strings like "1.000" or "6.000" would
#  otherwise be auto-converted to numbers by Excel and lose trailing zeros).
function Set-CellText {
    param($cell, [string]$text)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-CellText $ws.Range("D2") "29.622.42"
Set-CellText $ws.Range("E2") "  +5.53%  "
Set-CellText $ws.Range("D3") "1.920.52"
Set-CellText $ws.Range("E3") "  +4.23%  "
Set-CellText $ws.Range("E4") "  -0.41%  "
Set-CellText $ws.Range("D5") "335.57"
Set-CellText $ws.Range("E5") "  +1.65%  "
Set-CellText $ws.Range("D6") "1.001"
Set-CellText $ws.Range("E6") "  -0.27%  "
Set-CellText $ws.Range("D7") "0.4678"
Set-CellText $ws.Range("E7") "  +3.41%  "
Set-CellText $ws.Range("D8") "0.4120"
Set-CellText $ws.Range("E8") "  +5.83%  "
Set-CellText $ws.Range("D9") "48.18"
Set-CellText $ws.Range("E9") "  +0.83%  "
Set-CellText $ws.Range("D10") "0.08036"
Set-CellText $ws.Range("E10") "  +3.68%  "
Set-CellText $ws.Range("D11") "1.014"
Set-CellText $ws.Range("E11") "  +4.47%  "
Set-CellText $ws.Range("D12") "22.42"
Set-CellText $ws.Range("E12") "  +5.61%  "
Set-CellText $ws.Range("D13") "1.928.84"
Set-CellText $ws.Range("E13") "  +4.72%  "
Set-CellText $ws.Range("D14") "6.000"
Set-CellText $ws.Range("E14") "  +4.31%  "
Set-CellText $ws.Range("D15") "7.183"
Set-CellText $ws.Range("E15") "  +3.72%  "
Set-CellText $ws.Range("D16") "89.96"
Set-CellText $ws.Range("E16") "  +3.66%  "
Set-CellText $ws.Range("E17") "  -0.26%  "
Set-CellText $ws.Range("D18") "0.00001035"
Set-CellText $ws.Range("E18") "  +2.18%  "
Set-CellText $ws.Range("D19") "0.06608"
Set-CellText $ws.Range("E19") "  +1.22%  "
Set-CellText $ws.Range("E20") "  +6.00%  "
Set-CellText $ws.Range("D21") "1.000"
Set-CellText $ws.Range("E21") "  -1.64%  "
Set-CellText $ws.Range("D22") "29.607.49"
Set-CellText $ws.Range("E22") "  +5.41%  "
Set-CellText $ws.Range("D23") "5.571"
Set-CellText $ws.Range("E23") "  +6.22%  "
Set-CellText $ws.Range("D24") "11.67"
Set-CellText $ws.Range("E24") "  +10.68%  "
Set-CellText $ws.Range("D25") "2.206"
Set-CellText $ws.Range("E25") "  -1.57%  "
Set-CellText $ws.Range("D26") "2.168.90"
Set-CellText $ws.Range("E26") "  +5.05%  "
Set-CellText $ws.Range("D27") "156.21"
Set-CellText $ws.Range("E27") "  +0.09%  "
Set-CellText $ws.Range("D28") "19.91"
Set-CellText $ws.Range("E28") "  +4.61%  "
Set-CellText $ws.Range("D29") "2.146"
Set-CellText $ws.Range("E29") "  +6.47%  "
Set-CellText $ws.Range("D30") "5.724"
Set-CellText $ws.Range("E30") "  +9.70%  "
Set-CellText $ws.Range("D31") "117.62"
Set-CellText $ws.Range("E31") "  +1.58%  "
Set-CellText $ws.Range("D32") "1.073"
Set-CellText $ws.Range("E32") "  +15.69%  "
Set-CellText $ws.Range("D33") "0.09492"
Set-CellText $ws.Range("E33") "  +3.13%  "
Set-CellText $ws.Range("D34") "1.437"
Set-CellText $ws.Range("E34") "  +5.31%  "
Set-CellText $ws.Range("D35") "3.574"
Set-CellText $ws.Range("E35") "  -0.71%  "
Set-CellText $ws.Range("D36") "5.418"
Set-CellText $ws.Range("E36") "  +5.32%  "
Set-CellText $ws.Range("D37") "0.06145"
Set-CellText $ws.Range("E37") "  +2.84%  "
Set-CellText $ws.Range("D38") "0.02276"
Set-CellText $ws.Range("E38") "  +4.51%  "
Set-CellText $ws.Range("D39") "8.431"
Set-CellText $ws.Range("E39") "  +4.48%  "
Set-CellText $ws.Range("D40") "1.181"
Set-CellText $ws.Range("E40") "  +1.29%  "
Set-CellText $ws.Range("D41") "0.5900"
Set-CellText $ws.Range("E41") "  +5.04%  "
Set-CellText $ws.Range("B42") "Algorand"
Set-CellText $ws.Range("C42") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-CellText $ws.Range("D42") "0.1849"
Set-CellText $ws.Range("E42") "  +4.20%  "
Set-CellText $ws.Range("B43") "Aptos"
Set-CellText $ws.Range("C43") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-CellText $ws.Range("D43") "10.23"
Set-CellText $ws.Range("E43") "  +3.36%  "
Set-CellText $ws.Range("B44") "WEMIXTOKEN"
Set-CellText $ws.Range("C44") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-CellText $ws.Range("D44") "1.261"
Set-CellText $ws.Range("E44") "  +1.83%  "
Set-CellText $ws.Range("B45") "RenderToken"
Set-CellText $ws.Range("C45") "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws.Range("D45") "2.349"
Set-CellText $ws.Range("E45") "  +4.10%  "
Set-CellText $ws.Range("B46") "Cronos"
Set-CellText $ws.Range("C46") "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText $ws.Range("D46") "0.07506"
Set-CellText $ws.Range("E46") "  +4.60%  "
Set-CellText $ws.Range("B47") "Decentraland"
Set-CellText $ws.Range("C47") "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-CellText $ws.Range("D47") "0.5586"
Set-CellText $ws.Range("E47") "  +4.82%  "
Set-CellText $ws.Range("B48") "EnergySwap"
Set-CellText $ws.Range("C48") "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws.Range("D48") "12.25"
Set-CellText $ws.Range("E48") "  +5.31%  "
Set-CellText $ws.Range("B49") "NEARProtocol"
Set-CellText $ws.Range("C49") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-CellText $ws.Range("D49") "1.936"
Set-CellText $ws.Range("E49") "  +4.30%  "
Set-CellText $ws.Range("B50") "Quant"
Set-CellText $ws.Range("C50") "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
Set-CellText $ws.Range("D50") "113.17"
Set-CellText $ws.Range("E50") "  +3.95%  "
Set-CellText $ws.Range("B51") "WOONetwork"
Set-CellText $ws.Range("C51") "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
Set-CellText $ws.Range("D51") "0.2995"
Set-CellText $ws.Range("E51") "  +14.64%  "
